# Update "Förändrad" (column C) dates from 2023-09-06 (serial 45175) to
# 2023-09-14 (serial 45183) for rows 2-5 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C5").Value = "2023-09-14"
